# Added ifo GDP component analysis preprocessing:
# a new year-over-year series column (CD) is appended after the existing
# CC column, extending the dimension from A1:CC36 to A1:CD36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CD1: new header date (2025-11-25, serial 45986). Copy CC1's formatting
# (bold/border/date-number-format style) first so CD1 shares the same
# cell style as the rest of row 1, then write the date value.
$ws.Range("CC1").Copy()
$ws.Range("CD1").PasteSpecial(-4122)
$ws.Range("CD1").Value = 45986

# CD4:CD35: new year-over-year growth figures for the added column.
# (Rows 2, 3 and 36 only hold the date in column A and have no
# corresponding CC/CD data point, so they are left untouched.)
$ws.Range("CD4").Value = -1.214503843588766
$ws.Range("CD5").Value = 2.320292790391942
$ws.Range("CD6").Value = 1.885678566467552
$ws.Range("CD7").Value = 0.8418828180919435
$ws.Range("CD8").Value = 2.156425336307732
$ws.Range("CD9").Value = 2.008337744838529
$ws.Range("CD10").Value = 1.767701339560834
$ws.Range("CD11").Value = 3.356068332483475
$ws.Range("CD12").Value = 2.096888587375512
$ws.Range("CD13").Value = 0.0906098039807901
$ws.Range("CD14").Value = -0.5371314765080459
$ws.Range("CD15").Value = 1.043657656017705
$ws.Range("CD16").Value = 1.074913676625111
$ws.Range("CD17").Value = 4.14401568487659
$ws.Range("CD18").Value = 3.462833019567579
$ws.Range("CD19").Value = 0.6671812968680912
$ws.Range("CD20").Value = -6.164758722681718
$ws.Range("CD21").Value = 4.439190964013684
$ws.Range("CD22").Value = 3.831321260898735
$ws.Range("CD23").Value = 0.7051540842417214
$ws.Range("CD24").Value = 0.6040035278082057
$ws.Range("CD25").Value = 2.231490582578455
$ws.Range("CD26").Value = 1.166748954083641
$ws.Range("CD27").Value = 2.21154834069659
$ws.Range("CD28").Value = 3.234237177393018
$ws.Range("CD29").Value = 1.234177215189858
$ws.Range("CD30").Value = 0.7261884723591683
$ws.Range("CD31").Value = -4.719616128339188
$ws.Range("CD32").Value = 4.008819402685915
$ws.Range("CD33").Value = 2.252360763152805
$ws.Range("CD34").Value = -0.313331919805826
$ws.Range("CD35").Value = -0.5482818925178212
